$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
Write-Output ("Slide 13 shapes: " + $s.Shapes.Count)
